$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the previously-empty Modularity table (Tabela3) data cells S5:V6
$ws.Range("S5").Value = 0.141
$ws.Range("T5").Value = 0.14
$ws.Range("U5").Value = "~0"
$ws.Range("V5").Value = "~0"

$ws.Range("S6").Value = 0.82
$ws.Range("U6").Value = "-"
$ws.Range("V6").Value = "~0"

# Apply the "0.000" number format to the full Modularity table data range
$ws.Range("S3:V6").NumberFormat = "0.000"

# 2. Add the new Accuracy table data (R9:V10) and turn it into Tabela4
$ws.Range("R9").Value = "Accuracy"
$ws.Range("S9").Value = "Louvain"
$ws.Range("T9").Value = "Leiden"
$ws.Range("U9").Value = "Girvan-Newman"
$ws.Range("V9").Value = "Infomap"
$ws.Range("R10").Value = "Email"
$ws.Range("S10").Value = 0.537
$ws.Range("T10").Value = 0.579
$ws.Range("U10").Value = 0.075
$ws.Range("V10").Value = 0.537

$rng4 = $ws.Range("R9:V10")
$tbl4 = $ws.ListObjects.Add(1, $rng4, $null, 1)
$tbl4.Name = "Tabela4"
$tbl4.TableStyle = "TableStyleLight13"

# 3. Add legend / reference text in column J
$ws.Range("J8").Value = "ACS - Average Community Size"
$ws.Range("J9").Value = "APL - Communities Average Path Length"
$ws.Range("J10").Value = "AID - Average Internal Degree"
$ws.Range("J11").Value = "AE - Average Embeddedness"
$ws.Range("J12").Value = "AT - Average Transitivity"

# 4. Minor cosmetic tweaks: column widths for S:T and selection
$ws.Range("S1:T1").ColumnWidth = 8.6
$ws.Range("I9").Borders.LineStyle = -4142
$ws.Range("Q15").Select()
